$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1976620616365569
$ws.Range("C2").Value = 0.5696068012752391
$ws.Range("J2").Value = 0.01806588735387885
$ws.Range("P2").Value = 0.1487778958554729
$ws.Range("S2").Value = 0.06588735387885228
$ws.Range("B3").Value = 0.00539568345323741
$ws.Range("C3").Value = 0.03057553956834532
$ws.Range("J3").Value = 0.0341726618705036
$ws.Range("P3").Value = 0.7535971223021583
$ws.Range("S3").Value = 0.1762589928057554
$ws.Range("J4").Value = 0.02013422818791946
$ws.Range("P4").Value = 0.7181208053691275
$ws.Range("S4").Value = 0.261744966442953
$ws.Range("B6").Value = 0.05900151285930409
$ws.Range("D6").Value = 0.01210287443267776
$ws.Range("E6").Value = 0.00151285930408472
$ws.Range("F6").Value = 0.06656580937972768
$ws.Range("J6").Value = 0.2753403933434191
$ws.Range("O6").Value = 0.01059001512859304
$ws.Range("Q6").Value = 0.1603630862329803
$ws.Range("R6").Value = 0.0529500756429652
$ws.Range("S6").Value = 0.3615733736762481
$ws.Range("B7").Value = 0.1270627062706271
$ws.Range("D7").Value = 0.0165016501650165
$ws.Range("E7").Value = 0.00165016501650165
$ws.Range("F7").Value = 0.0396039603960396
$ws.Range("J7").Value = 0.1303630363036304
$ws.Range("O7").Value = 0.0132013201320132
$ws.Range("Q7").Value = 0.1831683168316832
$ws.Range("R7").Value = 0.08415841584158416
$ws.Range("S7").Value = 0.4042904290429043
$ws.Range("B8").Value = 0.1042760690172543
$ws.Range("D8").Value = 0.02475618904726181
$ws.Range("E8").Value = 0.002250562640660165
$ws.Range("F8").Value = 0.06001500375093773
$ws.Range("J8").Value = 0.09302325581395349
$ws.Range("O8").Value = 0.01800450112528132
$ws.Range("Q8").Value = 0.1747936984246062
$ws.Range("R8").Value = 0.08927231807951988
$ws.Range("S8").Value = 0.4336084021005251
$ws.Range("B9").Value = 0.1099830795262267
$ws.Range("D9").Value = 0.01353637901861252
$ws.Range("F9").Value = 0.05922165820642978
$ws.Range("J9").Value = 0.1015228426395939
$ws.Range("O9").Value = 0.02030456852791878
$ws.Range("Q9").Value = 0.1878172588832487
$ws.Range("R9").Value = 0.08629441624365482
$ws.Range("S9").Value = 0.4213197969543147
$ws.Range("B10").Value = 0.1102040816326531
$ws.Range("D10").Value = 0.02295918367346939
$ws.Range("E10").Value = 0.001020408163265306
$ws.Range("F10").Value = 0.06428571428571428
$ws.Range("J10").Value = 0.1038265306122449
$ws.Range("O10").Value = 0.01709183673469388
$ws.Range("Q10").Value = 0.2188775510204082
$ws.Range("R10").Value = 0.08801020408163265
$ws.Range("S10").Value = 0.3737244897959184
$ws.Range("G11").Value = 0.1417830290010741
$ws.Range("J11").Value = 0.07841031149301826
$ws.Range("K11").Value = 0.1901181525241676
$ws.Range("L11").Value = 0.5757250268528464
$ws.Range("S11").Value = 0.01396348012889366
$ws.Range("G12").Value = 0.7644927536231884
$ws.Range("J12").Value = 0.1793478260869565
$ws.Range("K12").Value = 0.007246376811594203
$ws.Range("L12").Value = 0.03442028985507246
$ws.Range("S12").Value = 0.01449275362318841
$ws.Range("G13").Value = 0.5897435897435898
$ws.Range("J13").Value = 0.358974358974359
$ws.Range("S13").Value = 0.05128205128205128
$ws.Range("F15").Value = 0.01940298507462686
$ws.Range("H15").Value = 0.1208955223880597
$ws.Range("I15").Value = 0.08208955223880597
$ws.Range("J15").Value = 0.3656716417910448
$ws.Range("K15").Value = 0.06865671641791045
$ws.Range("M15").Value = 0.01194029850746269
$ws.Range("N15").Value = 0.001492537313432836
$ws.Range("O15").Value = 0.07611940298507462
$ws.Range("S15").Value = 0.2537313432835821
$ws.Range("F16").Value = 0.02764976958525346
$ws.Range("H16").Value = 0.17357910906298
$ws.Range("I16").Value = 0.08448540706605223
$ws.Range("J16").Value = 0.3963133640552995
$ws.Range("K16").Value = 0.1075268817204301
$ws.Range("M16").Value = 0.01536098310291859
$ws.Range("O16").Value = 0.06451612903225806
$ws.Range("S16").Value = 0.130568356374808
$ws.Range("F17").Value = 0.02677942212825934
$ws.Range("H17").Value = 0.1592670894996476
$ws.Range("I17").Value = 0.08668076109936575
$ws.Range("J17").Value = 0.4418604651162791
$ws.Range("K17").Value = 0.09795630725863284
$ws.Range("M17").Value = 0.0204369274136716
$ws.Range("O17").Value = 0.05637773079633545
$ws.Range("S17").Value = 0.1106412966878083
$ws.Range("F18").Value = 0.01848739495798319
$ws.Range("H18").Value = 0.134453781512605
$ws.Range("I18").Value = 0.09411764705882353
$ws.Range("J18").Value = 0.4436974789915966
$ws.Range("K18").Value = 0.09747899159663866
$ws.Range("M18").Value = 0.01008403361344538
$ws.Range("N18").Value = 0.001680672268907563
$ws.Range("O18").Value = 0.06218487394957983
$ws.Range("S18").Value = 0.1378151260504202
$ws.Range("F19").Value = 0.01559703400664792
$ws.Range("H19").Value = 0.2155458961902327
$ws.Range("I19").Value = 0.07951930452569675
$ws.Range("J19").Value = 0.3750958834057786
$ws.Range("K19").Value = 0.108667859882383
$ws.Range("M19").Value = 0.01789823574533368
$ws.Range("N19").Value = 0.0005113781641523907
$ws.Range("O19").Value = 0.06775760675019177
$ws.Range("S19").Value = 0.1194068013295832
